$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1288654.33
$ws.Range("B3").Value = 354787.37
$ws.Range("B4").Value = 19882.9
$ws.Range("B5").Value = 6533.7
$ws.Range("B6").Value = 1401.27
$ws.Range("B7").Value = 281.83
